$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '23.185.07'
$c.NumberFormat = "General"

$c = $ws.Range("E2")
$c.NumberFormat = "@"
$c.Value = '  -1.13%  '
$c.NumberFormat = "General"

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '1.612.07'
$c.NumberFormat = "General"

$c = $ws.Range("E3")
$c.NumberFormat = "@"
$c.Value = '  -1.18%  '
$c.NumberFormat = "General"

$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.NumberFormat = "General"

$c = $ws.Range("E4")
$c.NumberFormat = "@"
$c.Value = '  +0.39%  '
$c.NumberFormat = "General"

$c = $ws.Range("E5")
$c.NumberFormat = "@"
$c.Value = '  +0.30%  '
$c.NumberFormat = "General"

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '303.01'
$c.NumberFormat = "General"

$c = $ws.Range("E6")
$c.NumberFormat = "@"
$c.Value = '  -0.55%  '
$c.NumberFormat = "General"

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '0.3783'
$c.NumberFormat = "General"

$c = $ws.Range("E7")
$c.NumberFormat = "@"
$c.Value = '  -0.28%  '
$c.NumberFormat = "General"

$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = '52.09'
$c.NumberFormat = "General"

$c = $ws.Range("E8")
$c.NumberFormat = "@"
$c.Value = '  +0.10%  '
$c.NumberFormat = "General"

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.3532'
$c.NumberFormat = "General"

$c = $ws.Range("E9")
$c.NumberFormat = "@"
$c.Value = '  -2.97%  '
$c.NumberFormat = "General"

$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = '0.08098'
$c.NumberFormat = "General"

$c = $ws.Range("E10")
$c.NumberFormat = "@"
$c.Value = '  -0.22%  '
$c.NumberFormat = "General"

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '1.204'
$c.NumberFormat = "General"

$c = $ws.Range("E11")
$c.NumberFormat = "@"
$c.Value = '  -2.37%  '
$c.NumberFormat = "General"

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '1.003'
$c.NumberFormat = "General"

$c = $ws.Range("E12")
$c.NumberFormat = "@"
$c.Value = '  +0.42%  '
$c.NumberFormat = "General"

$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = '22.11'
$c.NumberFormat = "General"

$c = $ws.Range("E13")
$c.NumberFormat = "@"
$c.Value = '  -2.63%  '
$c.NumberFormat = "General"

$c = $ws.Range("E14")
$c.NumberFormat = "@"
$c.Value = '  -3.36%  '
$c.NumberFormat = "General"

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '7.255'
$c.NumberFormat = "General"

$c = $ws.Range("E15")
$c.NumberFormat = "@"
$c.Value = '  +0.13%  '
$c.NumberFormat = "General"

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '0.00001209'
$c.NumberFormat = "General"

$c = $ws.Range("E16")
$c.NumberFormat = "@"
$c.Value = '  -3.25%  '
$c.NumberFormat = "General"

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '1.598.99'
$c.NumberFormat = "General"

$c = $ws.Range("E17")
$c.NumberFormat = "@"
$c.Value = '  -1.42%  '
$c.NumberFormat = "General"

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '94.31'
$c.NumberFormat = "General"

$c = $ws.Range("E18")
$c.NumberFormat = "@"
$c.Value = '  +0.57%  '
$c.NumberFormat = "General"

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '0.06908'
$c.NumberFormat = "General"

$c = $ws.Range("E19")
$c.NumberFormat = "@"
$c.Value = '  -0.43%  '
$c.NumberFormat = "General"

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '6.491'
$c.NumberFormat = "General"

$c = $ws.Range("E20")
$c.NumberFormat = "@"
$c.Value = '  +0.80%  '
$c.NumberFormat = "General"

$ws.Range("B21").Value = 'Dai'

$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'

$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = '1.002'
$c.NumberFormat = "General"

$c = $ws.Range("E21")
$c.NumberFormat = "@"
$c.Value = '  +0.36%  '
$c.NumberFormat = "General"

$ws.Range("B22").Value = 'Avalanche'

$ws.Range("C22").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'

$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = '17.19'
$c.NumberFormat = "General"

$c = $ws.Range("E22")
$c.NumberFormat = "@"
$c.Value = '  -4.21%  '
$c.NumberFormat = "General"

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '12.30'
$c.NumberFormat = "General"

$c = $ws.Range("E23")
$c.NumberFormat = "@"
$c.Value = '  -3.59%  '
$c.NumberFormat = "General"

$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = '23.173.80'
$c.NumberFormat = "General"

$c = $ws.Range("E24")
$c.NumberFormat = "@"
$c.Value = '  -1.20%  '
$c.NumberFormat = "General"

$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = '2.506'
$c.NumberFormat = "General"

$c = $ws.Range("E25")
$c.NumberFormat = "@"
$c.Value = '  +3.71%  '
$c.NumberFormat = "General"

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '3.007'
$c.NumberFormat = "General"

$c = $ws.Range("E26")
$c.NumberFormat = "@"
$c.Value = '  -6.72%  '
$c.NumberFormat = "General"

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '20.87'
$c.NumberFormat = "General"

$c = $ws.Range("E27")
$c.NumberFormat = "@"
$c.Value = '  -1.65%  '
$c.NumberFormat = "General"

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '150.88'
$c.NumberFormat = "General"

$c = $ws.Range("E28")
$c.NumberFormat = "@"
$c.Value = '  +0.77%  '
$c.NumberFormat = "General"

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '5.232'
$c.NumberFormat = "General"

$c = $ws.Range("E29")
$c.NumberFormat = "@"
$c.Value = '  -0.94%  '
$c.NumberFormat = "General"

$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = '132.21'
$c.NumberFormat = "General"

$c = $ws.Range("E30")
$c.NumberFormat = "@"
$c.Value = '  -1.96%  '
$c.NumberFormat = "General"

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '1.782.36'
$c.NumberFormat = "General"

$c = $ws.Range("E31")
$c.NumberFormat = "@"
$c.Value = '  -1.03%  '
$c.NumberFormat = "General"

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '1.066'
$c.NumberFormat = "General"

$c = $ws.Range("E32")
$c.NumberFormat = "@"
$c.Value = '  +11.07%  '
$c.NumberFormat = "General"

$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = '6.464'
$c.NumberFormat = "General"

$c = $ws.Range("E33")
$c.NumberFormat = "@"
$c.Value = '  -5.43%  '
$c.NumberFormat = "General"

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '2.100'
$c.NumberFormat = "General"

$c = $ws.Range("E34")
$c.NumberFormat = "@"
$c.Value = '  -8.82%  '
$c.NumberFormat = "General"

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '11.29'
$c.NumberFormat = "General"

$c = $ws.Range("E35")
$c.NumberFormat = "@"
$c.Value = '  +1.82%  '
$c.NumberFormat = "General"

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '0.02701'
$c.NumberFormat = "General"

$c = $ws.Range("E36")
$c.NumberFormat = "@"
$c.Value = '  -3.66%  '
$c.NumberFormat = "General"

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.08741'
$c.NumberFormat = "General"

$c = $ws.Range("E37")
$c.NumberFormat = "@"
$c.Value = '  -1.29%  '
$c.NumberFormat = "General"

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.2450'
$c.NumberFormat = "General"

$c = $ws.Range("E38")
$c.NumberFormat = "@"
$c.Value = '  -3.52%  '
$c.NumberFormat = "General"

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '0.06922'
$c.NumberFormat = "General"

$c = $ws.Range("E39")
$c.NumberFormat = "@"
$c.Value = '  -4.03%  '
$c.NumberFormat = "General"

$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = '5.839'
$c.NumberFormat = "General"

$c = $ws.Range("E40")
$c.NumberFormat = "@"
$c.Value = '  -4.73%  '
$c.NumberFormat = "General"

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '1.322'
$c.NumberFormat = "General"

$c = $ws.Range("E41")
$c.NumberFormat = "@"
$c.Value = '  -2.84%  '
$c.NumberFormat = "General"

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '0.6876'
$c.NumberFormat = "General"

$c = $ws.Range("E42")
$c.NumberFormat = "@"
$c.Value = '  -3.25%  '
$c.NumberFormat = "General"

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '11.93'
$c.NumberFormat = "General"

$c = $ws.Range("E43")
$c.NumberFormat = "@"
$c.Value = '  -3.60%  '
$c.NumberFormat = "General"

$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = '15.27'
$c.NumberFormat = "General"

$c = $ws.Range("E44")
$c.NumberFormat = "@"
$c.Value = '  -6.50%  '
$c.NumberFormat = "General"

$c = $ws.Range("E45")
$c.NumberFormat = "@"
$c.Value = '  +0.27%  '
$c.NumberFormat = "General"

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '0.6308'
$c.NumberFormat = "General"

$c = $ws.Range("E46")
$c.NumberFormat = "@"
$c.Value = '  -3.25%  '
$c.NumberFormat = "General"

$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = '3.945'
$c.NumberFormat = "General"

$c = $ws.Range("E47")
$c.NumberFormat = "@"
$c.Value = '  -1.37%  '
$c.NumberFormat = "General"

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '2.246'
$c.NumberFormat = "General"

$c = $ws.Range("E48")
$c.NumberFormat = "@"
$c.Value = '  -3.95%  '
$c.NumberFormat = "General"

$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = '0.07869'
$c.NumberFormat = "General"

$c = $ws.Range("E49")
$c.NumberFormat = "@"
$c.Value = '  -1.76%  '
$c.NumberFormat = "General"

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '126.44'
$c.NumberFormat = "General"

$c = $ws.Range("E50")
$c.NumberFormat = "@"
$c.Value = '  +0.44%  '
$c.NumberFormat = "General"

$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = '1.169'
$c.NumberFormat = "General"

$c = $ws.Range("E51")
$c.NumberFormat = "@"
$c.Value = '  -3.29%  '
$c.NumberFormat = "General"
